# Auto-generated edit script: updates the cryptocurrency price/volume
# table on Sheet1 to match the "Fri Mar 24 21:00:26 UTC 2023" GitHub
# Actions refresh (see commit message). Rows 46/47 (PancakeSwap /
# Decentraland) also swap places in the ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.738.63"
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("D3").Value = "1.757.73"
$ws.Range("E3").Value = "  -3.15%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "322.82"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").Value = "'0.4250"
$ws.Range("E7").Value = "  -3.81%  "
$ws.Range("D8").Value = "0.3631"
$ws.Range("E8").Value = "  -1.89%  "
$ws.Range("D9").Value = "0.07527"
$ws.Range("E9").Value = "  -2.40%  "
$ws.Range("D10").Value = "42.47"
$ws.Range("E10").Value = "  -5.28%  "
$ws.Range("D11").Value = "1.094"
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "20.63"
$ws.Range("E13").Value = "  -6.67%  "
$ws.Range("D14").Value = "'6.030"
$ws.Range("E14").Value = "  -3.50%  "
$ws.Range("D15").Value = "7.269"
$ws.Range("E15").Value = "  -3.92%  "
$ws.Range("D16").Value = "1.781.75"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("D17").Value = "91.24"
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D19").Value = "0.06378"
$ws.Range("E19").Value = "  -3.12%  "
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").Value = "17.03"
$ws.Range("E21").Value = "  -2.84%  "
$ws.Range("D22").Value = "5.892"
$ws.Range("E22").Value = "  -5.06%  "
$ws.Range("D23").Value = "27.777.00"
$ws.Range("E23").Value = "  -1.91%  "
$ws.Range("D24").Value = "11.18"
$ws.Range("E24").Value = "  -4.38%  "
$ws.Range("D25").Value = "2.098"
$ws.Range("E25").Value = "  +4.34%  "
$ws.Range("D26").Value = "160.37"
$ws.Range("E26").Value = "  +2.67%  "
$ws.Range("D27").Value = "20.24"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").Value = "1.974.66"
$ws.Range("E28").Value = "  -2.54%  "
$ws.Range("D29").Value = "2.132"
$ws.Range("E29").Value = "  -7.99%  "
$ws.Range("D30").Value = "124.85"
$ws.Range("E30").Value = "  -2.47%  "
$ws.Range("D31").Value = "'1.110"
$ws.Range("E31").Value = "  -7.70%  "
$ws.Range("D32").Value = "3.664"
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").Value = "5.563"
$ws.Range("E33").Value = "  -5.34%  "
$ws.Range("D34").Value = "0.08891"
$ws.Range("E34").Value = "  -3.65%  "
$ws.Range("D35").Value = "12.21"
$ws.Range("E35").Value = "  -6.62%  "
$ws.Range("D36").Value = "0.02298"
$ws.Range("E36").Value = "  -2.32%  "
$ws.Range("D37").Value = "0.2101"
$ws.Range("E37").Value = "  -3.18%  "
$ws.Range("D38").Value = "0.06028"
$ws.Range("E38").Value = "  -2.78%  "
$ws.Range("D39").Value = "0.6328"
$ws.Range("E39").Value = "  -3.91%  "
$ws.Range("D40").Value = "4.958"
$ws.Range("D41").Value = "1.185"
$ws.Range("E41").Value = "  -1.14%  "
$ws.Range("D42").Value = "1.003"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").Value = "7.907"
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("D44").Value = "1.396"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").Value = "13.35"
$ws.Range("E45").Value = "  -3.64%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.5872"
$ws.Range("E46").Value = "  -3.42%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.699"
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("D48").Value = "1.986"
$ws.Range("E48").Value = "  -2.47%  "
$ws.Range("D49").Value = "123.31"
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("D50").Value = "1.169"
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("D51").Value = "0.06826"
$ws.Range("E51").Value = "  -2.29%  "
